$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") bumps from 46065 to 46066 for all data rows.
$ws.Range("C2").Value = 46066
$ws.Range("C3").Value = 46066
$ws.Range("C4").Value = 46066
$ws.Range("C5").Value = 46066
$ws.Range("C6").Value = 46066
$ws.Range("C7").Value = 46066
$ws.Range("C8").Value = 46066
$ws.Range("C9").Value = 46066

# Rows 4-8 get re-ordered (identities shuffle positions) and the
# "Markägare" (F) value follows the identity it belongs to.
$ws.Range("A4").Value = "A 6983-2023"
$ws.Range("B4").Value = 44967.68585648148
$ws.Range("G4").Value = 5.4
$ws.Range("F4").ClearContents()

$ws.Range("A5").Value = "A 25254-2025"
$ws.Range("B5").Value = 45800.50479166667
$ws.Range("G5").Value = 0.2
$ws.Range("F5").Value = "Kommuner"

$ws.Range("A6").Value = "A 26074-2025"
$ws.Range("B6").Value = 45805.32366898148
$ws.Range("G6").Value = 1.3
$ws.Range("F6").ClearContents()

$ws.Range("A7").Value = "A 25251-2025"
$ws.Range("B7").Value = 45800.50082175926
$ws.Range("G7").Value = 0.7
$ws.Range("F7").Value = "Kommuner"

$ws.Range("A8").Value = "A 35734-2023"
$ws.Range("B8").Value = 45147.89258101852
$ws.Range("G8").Value = 5.9
$ws.Range("F8").ClearContents()
